# feat(import): include product note in import feature
#
# Insert a new "note" column between the existing "salesPrice" (I)
# and "isBarred" (J) columns on the product import example sheet.
# This shifts the old "isBarred" column (and its sample row value)
# one column to the right, and adds the new header "note" in its place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "isBarred" column (J),
# shifting isBarred (and the sample "nej" value below it) to column K.
$ws.Columns("J").Insert()

# Add the new header for the note column.
$ws.Range("J1").Value = "note"

# Reflect the cursor/selection position left behind after the edit.
$ws.Range("K7").Select()
